$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.406.56"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.18%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.676.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.43%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.02%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5304"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.23%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2704"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.78%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06401"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.53%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.04%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07806"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.02%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.506"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.70%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.673.79"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5556"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅8310"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.21%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.460.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.48%  "

$ws.Range("E18").Value = "  -0.01%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.734"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.90%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.63%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.64%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.339"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.01%  "

$ws.Range("E23").Value = "  +0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "142.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.47%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1283"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +5.43%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.409"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.14%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.69%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.422"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.06248"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.29%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.271"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.89%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.607"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.69%  "

$ws.Range("E32").Value = "  +1.08%  "

$ws.Range("E33").Value = "  +2.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.006"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.428"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.33%  "

$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6055"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.42%  "

$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.779"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.53%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.192"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +9.14%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01632"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.74%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.084.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.55%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8642"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.83%  "

$ws.Range("E42").Value = "  -0.06%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.44%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.821.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.96%  "

$ws.Range("E45").Value = "  +1.48%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "57.05"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.10%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.152"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.59%  "

$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.000"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.32%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05206"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.71%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.028"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.37%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4240"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.44%  "
